$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 used to be a 12-column header (Day..Course). The new layout only keeps
# a single title cell in A1; the rest of the old header row is dropped.
$ws.Range("A1").Value = "Herald College Kathmandu"
$ws.Range("B1:L1").ClearContents()

# Rows 2-10: the data columns get reshuffled.
#   old C (Module Code)  -> new D
#   old D (Module Title) -> new E
#   old E (Hours, number) -> new C
#   old H (Room)          -> new J
#   old J (Group)         -> new H
#   I (Block) stays put; F/G (Class Type/Lecturer) stay put
#   K (Level) and L (Course) are removed entirely
# Columns: A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8 I=9 J=10 K=11 L=12
for ($r = 2; $r -le 10; $r++) {
    $oldC = $ws.Cells.Item($r, 3).Value2
    $oldD = $ws.Cells.Item($r, 4).Value2
    $oldE = $ws.Cells.Item($r, 5).Value2
    $oldH = $ws.Cells.Item($r, 8).Value2
    $oldJ = $ws.Cells.Item($r, 10).Value2

    $ws.Cells.Item($r, 3).Value = $oldE
    $ws.Cells.Item($r, 4).Value = $oldC
    $ws.Cells.Item($r, 5).Value = $oldD
    $ws.Cells.Item($r, 8).Value = $oldJ
    $ws.Cells.Item($r, 10).Value = $oldH
}

# Drop the now-obsolete Level (K) and Course (L) columns entirely.
$ws.Range("K1:L10").ClearContents()
